$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the VIN value in A2 (was "7MSRP15H&V", becomes "7PRTL15H&V") ---
$ws.Range("A2").Value = "7PRTL15H&V"

# --- Capture old comment (anchored at AA2: STAT/COLL_SYMBOL area) before column deletion ---
$oldComment = $ws.Range("AA2").Comment()
$commentText = $oldComment.Text()
$oldComment.Delete()

# --- Delete column Z (STAT) entirely ---
$ws.Columns("Z").Delete()

# --- Delete column AB (CHOICE_TIER, was AC before the first deletion) entirely ---
$ws.Columns("AB").Delete()

# --- Re-create the comment at its new location Z2 (shifted left by the Z-column deletion) ---
$newComment = $ws.Range("Z2").AddComment($commentText)
$newComment.Shape.TextFrame.Characters(1, 21).Font.Bold = $true
$newComment.Shape.TextFrame.Characters(1, 21).Font.Name = "Tahoma"
$newComment.Shape.TextFrame.Characters(1, 21).Font.Size = 9
$newComment.Shape.TextFrame.Characters(22, 1000).Font.Name = "Tahoma"
$newComment.Shape.TextFrame.Characters(22, 1000).Font.Size = 9
$newComment.Visible = $false

# --- Update the view: scroll so column S is first visible, select Z2 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 19
$win.ScrollRow = 1
$ws.Range("Z2").Select()
